$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4:D6").Value = "изменено"
$ws.Range("D6").Select()
